$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows above the current row 681, shifting the existing
# rows 681:692 down to 685:696 (formats copied from the row above, matching
# Excel's default Insert behaviour - column D keeps the date style).
$ws.Rows("681:684").Insert()

# Row 681 (new): Lechuga / Escarola / Primera
$ws.Range("A681").Value = 1
$ws.Range("B681").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C681").Value = "Arica y Parinacota"
$ws.Range("D681").Value = 44628
$ws.Range("E681").Value = 15
$ws.Range("F681").Value = 100112033
$ws.Range("G681").Value = "Lechuga"
$ws.Range("H681").Value = "Escarola"
$ws.Range("I681").Value = "Primera"
$ws.Range("J681").Value = 120
$ws.Range("K681").Value = 4000
$ws.Range("L681").Value = 4500
$ws.Range("M681").Value = 4250
$ws.Range("N681").Value = "$/caja 12 unidades"
$ws.Range("O681").Value = "Región de Arica y Parinacota"
$ws.Range("P681").Value = 354
$ws.Range("Q681").Value = 12
$ws.Range("R681").Value = "Hortaliza"

# Row 682 (new): Lechuga / Escarola / Segunda
$ws.Range("A682").Value = 1
$ws.Range("B682").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C682").Value = "Arica y Parinacota"
$ws.Range("D682").Value = 44628
$ws.Range("E682").Value = 15
$ws.Range("F682").Value = 100112033
$ws.Range("G682").Value = "Lechuga"
$ws.Range("H682").Value = "Escarola"
$ws.Range("I682").Value = "Segunda"
$ws.Range("J682").Value = 130
$ws.Range("K682").Value = 4000
$ws.Range("L682").Value = 4500
$ws.Range("M682").Value = 4250
$ws.Range("N682").Value = "$/caja 18 unidades"
$ws.Range("O682").Value = "Región de Arica y Parinacota"
$ws.Range("P682").Value = 236
$ws.Range("Q682").Value = 18
$ws.Range("R682").Value = "Hortaliza"

# Row 683 (new): Lechuga / Marina / Primera
$ws.Range("A683").Value = 1
$ws.Range("B683").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C683").Value = "Arica y Parinacota"
$ws.Range("D683").Value = 44628
$ws.Range("E683").Value = 15
$ws.Range("F683").Value = 100112033
$ws.Range("G683").Value = "Lechuga"
$ws.Range("H683").Value = "Marina"
$ws.Range("I683").Value = "Primera"
$ws.Range("J683").Value = 130
$ws.Range("K683").Value = 5000
$ws.Range("L683").Value = 6000
$ws.Range("M683").Value = 5500
$ws.Range("N683").Value = "$/caja 12 unidades"
$ws.Range("O683").Value = "Región de Arica y Parinacota"
$ws.Range("P683").Value = 458
$ws.Range("Q683").Value = 12
$ws.Range("R683").Value = "Hortaliza"

# Row 684 (new): Lechuga / Marina / Segunda
$ws.Range("A684").Value = 1
$ws.Range("B684").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C684").Value = "Arica y Parinacota"
$ws.Range("D684").Value = 44628
$ws.Range("E684").Value = 15
$ws.Range("F684").Value = 100112033
$ws.Range("G684").Value = "Lechuga"
$ws.Range("H684").Value = "Marina"
$ws.Range("I684").Value = "Segunda"
$ws.Range("J684").Value = 120
$ws.Range("K684").Value = 5000
$ws.Range("L684").Value = 6000
$ws.Range("M684").Value = 5500
$ws.Range("N684").Value = "$/caja 18 unidades"
$ws.Range("O684").Value = "Región de Arica y Parinacota"
$ws.Range("P684").Value = 306
$ws.Range("Q684").Value = 18
$ws.Range("R684").Value = "Hortaliza"
